$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B3").Value = 'RGB'
$ws.Range("C3").Value = 'showimage'
$ws.Range("E3").Value = 128
$ws.Range("F3").Value = 3.0482
$ws.Range("G3").Value = 4
$ws.Range("A4").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B4").Value = 'RGB'
$ws.Range("C4").Value = 'leucolinf'
$ws.Range("E4").Value = 128
$ws.Range("F4").Value = 9.9551
$ws.Range("G4").Value = 2
$ws.Range("A5").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B5").Value = 'RGB'
$ws.Range("C5").Value = 'plain'
$ws.Range("E5").Value = 128
$ws.Range("F5").Value = 1.912
$ws.Range("G5").Value = 7
$ws.Range("A6").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B6").Value = 'RGB'
$ws.Range("C6").Value = 60079
$ws.Range("E6").Value = 128
$ws.Range("F6").Value = 7.3142
$ws.Range("G6").Value = 2
$ws.Range("A7").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B7").Value = 'RGB'
$ws.Range("C7").Value = 45096
$ws.Range("E7").Value = 128
$ws.Range("F7").Value = 2.272
$ws.Range("G7").Value = 5
$ws.Range("A8").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B8").Value = 'RGB'
$ws.Range("C8").Value = 188091
$ws.Range("E8").Value = 128
$ws.Range("F8").Value = 5.8066
$ws.Range("G8").Value = 4
$ws.Range("A9").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B9").Value = 'RGB'
$ws.Range("C9").Value = 'test'
$ws.Range("E9").Value = 128
$ws.Range("F9").Value = 3.2934
$ws.Range("G9").Value = 4
$ws.Range("A10").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B10").Value = 'RGB'
$ws.Range("C10").Value = 253036
$ws.Range("E10").Value = 128
$ws.Range("F10").Value = 4.5783
$ws.Range("G10").Value = 4
$ws.Range("A11").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B11").Value = 'RGB'
$ws.Range("C11").Value = 42049
$ws.Range("E11").Value = 128
$ws.Range("F11").Value = 4.3924
$ws.Range("G11").Value = 4
$ws.Range("A12").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B12").Value = 'RGB'
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 128
$ws.Range("A13").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B13").Value = 'RGB'
$ws.Range("C13").Value = 163014
$ws.Range("E13").Value = 128
$ws.Range("F13").Value = 3.8043
$ws.Range("G13").Value = 4
$ws.Range("A14").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B14").Value = 'RGB'
$ws.Range("C14").Value = 124084
$ws.Range("E14").Value = 128
$ws.Range("F14").Value = 2.0809
$ws.Range("G14").Value = 5
$ws.Range("A15").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B15").Value = 'RGB'
$ws.Range("C15").Value = 176035
$ws.Range("E15").Value = 128
$ws.Range("F15").Value = 3.0331
$ws.Range("G15").Value = 4
$ws.Range("A16").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B16").Value = 'RGB'
$ws.Range("C16").Value = 295087
$ws.Range("E16").Value = 128
$ws.Range("F16").Value = 4.1481
$ws.Range("G16").Value = 4
$ws.Range("A17").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B17").Value = 'RGB'
$ws.Range("C17").Value = 216066
$ws.Range("E17").Value = 128
$ws.Range("F17").Value = 4.5588
$ws.Range("G17").Value = 4
$ws.Range("A18").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B18").Value = 'RGB'
$ws.Range("C18").Value = '41004-2'
$ws.Range("E18").Value = 128
$ws.Range("F18").Value = 3.586
$ws.Range("G18").Value = 4
$ws.Range("A19").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B19").Value = 'RGB'
$ws.Range("C19").Value = 388016
$ws.Range("E19").Value = 128
$ws.Range("F19").Value = 4.1206
$ws.Range("G19").Value = 4
$ws.Range("A20").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B20").Value = 'RGB'
$ws.Range("C20").Value = 385028
$ws.Range("E20").Value = 128
$ws.Range("F20").Value = 4.2629
$ws.Range("G20").Value = 4
$ws.Range("A21").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B21").Value = 'RGB'
$ws.Range("C21").Value = 225017
$ws.Range("E21").Value = 128
$ws.Range("F21").Value = 3.8323
$ws.Range("G21").Value = 5
$ws.Range("A22").Value = 'metodo 1: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("B22").Value = 'RGB'
$ws.Range("C22").Value = 113044
$ws.Range("E22").Value = 128
$ws.Range("F22").Value = 6.3636
$ws.Range("G22").Value = 3
$ws.Range("A23").Value = 'fabric'
$ws.Range("B23").Value = 'RGB'
$ws.Range("C23").Value = 'metodo 2: k-means-Reagrupamiento-RGB-Mahalanobis'
$ws.Range("E23").Value = 64
$ws.Range("F23").Value = 2.1094
$ws.Range("G23").Value = 6
$ws.Range("C24").Value = '12375764_10154354426419428_2121622626_o'
$ws.Range("F24").Value = 2.8155
$ws.Range("G24").Value = 5
$ws.Range("C25").Value = 113044
$ws.Range("F25").Value = 3.3246
$ws.Range("G25").Value = 4
$ws.Range("C26").Value = 225017
$ws.Range("F26").Value = 2.0082
$ws.Range("G26").Value = 5
$ws.Range("C27").Value = 385028
$ws.Range("F27").Value = 2.2036
$ws.Range("G27").Value = 6
$ws.Range("C28").Value = 388016
$ws.Range("F28").Value = 2.3084
$ws.Range("G28").Value = 4
$ws.Range("C29").Value = '41004-2'
$ws.Range("F29").Value = 1.8198
$ws.Range("G29").Value = 6

# Update view state: scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("D35").Select()
